$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 3799.182
$ws.Range("I15").Value2 = 3799.182
$ws.Range("K15").Value2 = 11397.546
$ws.Range("M15").Value2 = -11228.546

$ws.Range("H32").Value2 = 965.6667
$ws.Range("I32").Value2 = 933.3333
$ws.Range("J32").Value2 = 998
$ws.Range("K32").Value2 = 933.3333
$ws.Range("L32").Value2 = 998
$ws.Range("M32").Value2 = -607.3333
$ws.Range("N32").Value2 = -1650

$ws.Range("H33").Value2 = 376.5238
$ws.Range("I33").Value2 = 311.05264
$ws.Range("K33").Value2 = 311.05264
$ws.Range("M33").Value2 = -82.05264

$ws.Range("H92").Value2 = 676.625
$ws.Range("I92").Value2 = 676.625
$ws.Range("K92").Value2 = 676.625
$ws.Range("M92").Value2 = 571.375

$ws.Range("H107").Value2 = 853.06665
$ws.Range("I107").Value2 = 755.8461
$ws.Range("J107").Value2 = 1485
$ws.Range("K107").Value2 = 755.8461
$ws.Range("L107").Value2 = 1485
$ws.Range("M107").Value2 = 1164.1539
$ws.Range("N107").Value2 = -5325

$ws.Range("H112").Value2 = 2832
$ws.Range("I112").Value2 = 0
$ws.Range("J112").Value2 = 2832
$ws.Range("K112").Value2 = 0
$ws.Range("L112").Value2 = 8496
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value2 = -10712

$ws.Range("H132").Value2 = 25284.5
$ws.Range("I132").Value2 = 3821.1904
$ws.Range("K132").Value2 = 11463.5712
$ws.Range("M132").Value2 = -8933.5712

$ws.Range("H137").Value2 = 12447.385
$ws.Range("I137").Value2 = 2388.3635
$ws.Range("J137").Value2 = 19824
$ws.Range("K137").Value2 = 7165.0905
$ws.Range("L137").Value2 = 59472
$ws.Range("M137").Value2 = -4615.0905
$ws.Range("N137").Value2 = -64572

$ws.Range("H138").Value2 = 5410.306
$ws.Range("I138").Value2 = 6664.5454
$ws.Range("J138").Value2 = 5047.237
$ws.Range("K138").Value2 = 19993.6362
$ws.Range("L138").Value2 = 15141.711
$ws.Range("M138").Value2 = -14853.6362
$ws.Range("N138").Value2 = -25421.711

$ws.Range("H141").Value2 = 3408.4546
$ws.Range("I141").Value2 = 3265.3333
$ws.Range("J141").Value2 = 4052.5
$ws.Range("K141").Value2 = 9795.999899999999
$ws.Range("L141").Value2 = 12157.5
$ws.Range("M141").Value2 = -4615.999899999999
$ws.Range("N141").Value2 = -22517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 123.2
$ws.Range("I5").Value2 = 109.111115
$ws.Range("J5").Value2 = 250
$ws.Range("K5").Value2 = 109.111115
$ws.Range("L5").Value2 = 250
$ws.Range("M5").Value2 = 2.888885000000002
$ws.Range("N5").Value2 = -474

$ws.Range("H32").Value2 = 1436479.4
$ws.Range("I32").Value2 = 640228.25
$ws.Range("K32").Value2 = 640228.25
$ws.Range("M32").Value2 = -639941.25

$ws.Range("H97").Value2 = 508.125
$ws.Range("I97").Value2 = 521.1539
$ws.Range("J97").Value2 = 451.66666
$ws.Range("K97").Value2 = 521.1539
$ws.Range("L97").Value2 = 451.66666
$ws.Range("M97").Value2 = -25.15390000000002
$ws.Range("N97").Value2 = -1443.66666

$ws.Range("H102").Value2 = 50001850
$ws.Range("I102").Value2 = 55557084
$ws.Range("K102").Value2 = 55557084
$ws.Range("M102").Value2 = -55555462

$ws.Range("H132").Value2 = 3205.05
$ws.Range("I132").Value2 = 2577.2
$ws.Range("K132").Value2 = 7731.599999999999
$ws.Range("M132").Value2 = -5201.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 123.2
$ws.Range("I4").Value2 = 109.111115
$ws.Range("J4").Value2 = 250
$ws.Range("K4").Value2 = 109.111115
$ws.Range("L4").Value2 = 250
$ws.Range("M4").Value2 = 5.888885000000002
$ws.Range("N4").Value2 = -480

$ws.Range("H94").Value2 = 111112320
$ws.Range("I94").Value2 = 121213350
$ws.Range("J94").Value2 = 1000
$ws.Range("K94").Value2 = 121213350
$ws.Range("L94").Value2 = 1000
$ws.Range("M94").Value2 = -121212899
$ws.Range("N94").Value2 = -1902

$ws.Range("H107").Value2 = 2482809.8
$ws.Range("I107").Value2 = 2959839.2
$ws.Range("K107").Value2 = 2959839.2
$ws.Range("M107").Value2 = -2957919.2

$ws.Range("H134").Value2 = 2740.4443
$ws.Range("I134").Value2 = 2770.5
$ws.Range("J134").Value2 = 2500
$ws.Range("K134").Value2 = 8311.5
$ws.Range("L134").Value2 = 7500
$ws.Range("M134").Value2 = -5776.5
$ws.Range("N134").Value2 = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2276616
$ws.Range("I31").Value2 = 2574.125
$ws.Range("J31").Value2 = 3209556.2
$ws.Range("K31").Value2 = 2574.125
$ws.Range("L31").Value2 = 3209556.2
$ws.Range("M31").Value2 = -2279.125
$ws.Range("N31").Value2 = -3210146.2

$ws.Range("H34").Value2 = 2276616
$ws.Range("I34").Value2 = 2574.125
$ws.Range("J34").Value2 = 3209556.2
$ws.Range("K34").Value2 = 2574.125
$ws.Range("L34").Value2 = 3209556.2
$ws.Range("M34").Value2 = -2372.125
$ws.Range("N34").Value2 = -3209960.2

$ws.Range("H58").Value2 = 2811.182
$ws.Range("I58").Value2 = 1989
$ws.Range("J58").Value2 = 4250
$ws.Range("K58").Value2 = 1989
$ws.Range("L58").Value2 = 4250
$ws.Range("M58").Value2 = -1786
$ws.Range("N58").Value2 = -4656

$ws.Range("H107").Value2 = 799.65
$ws.Range("I107").Value2 = 612.1539
$ws.Range("J107").Value2 = 1147.8572
$ws.Range("K107").Value2 = 612.1539
$ws.Range("L107").Value2 = 1147.8572
$ws.Range("M107").Value2 = 1307.8461
$ws.Range("N107").Value2 = -4987.8572

$ws.Range("H134").Value2 = 3530.2424
$ws.Range("I134").Value2 = 3704.2917
$ws.Range("J134").Value2 = 3066.111
$ws.Range("K134").Value2 = 11112.8751
$ws.Range("L134").Value2 = 9198.332999999999
$ws.Range("M134").Value2 = -8577.875100000001
$ws.Range("N134").Value2 = -14268.333

$ws.Range("H136").Value2 = 2811.182
$ws.Range("I136").Value2 = 1989
$ws.Range("J136").Value2 = 4250
$ws.Range("K136").Value2 = 5967
$ws.Range("L136").Value2 = 12750
$ws.Range("M136").Value2 = -3417
$ws.Range("N136").Value2 = -17850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 1644647.8
$ws.Range("I68").Value2 = 4820.6665
$ws.Range("J68").Value2 = 2046238
$ws.Range("K68").Value2 = 14461.9995
$ws.Range("L68").Value2 = 6138714
$ws.Range("M68").Value2 = -13650.9995
$ws.Range("N68").Value2 = -6140336

$ws.Range("H71").Value2 = 1644647.8
$ws.Range("I71").Value2 = 4820.6665
$ws.Range("J71").Value2 = 2046238
$ws.Range("K71").Value2 = 43385.9985
$ws.Range("L71").Value2 = 18416142
$ws.Range("M71").Value2 = -39329.9985
$ws.Range("N71").Value2 = -18424254

$ws.Range("H118").Value2 = 12999.9
$ws.Range("I118").Value2 = 17276.428
$ws.Range("K118").Value2 = 51829.284
$ws.Range("M118").Value2 = -50586.284

$ws.Range("H131").Value2 = 2458784
$ws.Range("I131").Value2 = 12173.857
$ws.Range("J131").Value2 = 5884038.5
$ws.Range("K131").Value2 = 36521.571
$ws.Range("L131").Value2 = 17652115.5
$ws.Range("M131").Value2 = -31481.571
$ws.Range("N131").Value2 = -17662195.5

$ws.Range("H137").Value2 = 5869.4
$ws.Range("J137").Value2 = 6856.4287
$ws.Range("L137").Value2 = 20569.2861
$ws.Range("N137").Value2 = -30769.2861

$ws.Range("H140").Value2 = 11311.389
$ws.Range("I140").Value2 = 3372.5715
$ws.Range("J140").Value2 = 16363.363
$ws.Range("K140").Value2 = 10117.7145
$ws.Range("L140").Value2 = 49090.089
$ws.Range("M140").Value2 = -4937.7145
$ws.Range("N140").Value2 = -59450.089

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 29415226
$ws.Range("I122").Value2 = 3025.48
$ws.Range("K122").Value2 = 9076.440000000001
$ws.Range("M122").Value2 = -6626.440000000001

$ws.Range("H132").Value2 = 1703.8
$ws.Range("I132").Value2 = 1712.1666
$ws.Range("J132").Value2 = 1670.3334
$ws.Range("K132").Value2 = 5136.4998
$ws.Range("L132").Value2 = 5011.0002
$ws.Range("M132").Value2 = -2606.4998
$ws.Range("N132").Value2 = -10071.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 63109.523
$ws.Range("I40").Value2 = 91178.64
$ws.Range("K40").Value2 = 91178.64
$ws.Range("M40").Value2 = -91042.64

$ws.Range("H46").Value2 = 222
$ws.Range("I46").Value2 = 222
$ws.Range("J46").Value2 = 222
$ws.Range("K46").Value2 = 222
$ws.Range("L46").Value2 = 222
$ws.Range("M46").Value2 = -34
$ws.Range("N46").Value2 = -598

$ws.Range("H100").Value2 = 6758133
$ws.Range("I100").Value2 = 6758133
$ws.Range("J100").Value2 = 0
$ws.Range("K100").Value2 = 6758133
$ws.Range("L100").Value2 = 0
$ws.Range("M100").Value2 = -6757592
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value2 = 2940.1
$ws.Range("I122").Value2 = 3122.4285
$ws.Range("J122").Value2 = 2514.6667
$ws.Range("K122").Value2 = 9367.2855
$ws.Range("L122").Value2 = 7544.000100000001
$ws.Range("M122").Value2 = -6917.2855
$ws.Range("N122").Value2 = -12444.0001

$ws.Range("H136").Value2 = 4917.591
$ws.Range("I136").Value2 = 4069.5293
$ws.Range("J136").Value2 = 7801
$ws.Range("K136").Value2 = 12208.5879
$ws.Range("L136").Value2 = 23403
$ws.Range("M136").Value2 = -9658.5879
$ws.Range("N136").Value2 = -28503

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 12624.6
$ws.Range("I126").Value2 = 14530.75
$ws.Range("K126").Value2 = 43592.25
$ws.Range("M126").Value2 = -41122.25

$ws.Range("H131").Value2 = 113999.5
$ws.Range("J131").Value2 = 113999.5
$ws.Range("L131").Value2 = 113999.5
$ws.Range("N131").Value2 = -124079.5

$ws.Range("H136").Value2 = 6669969
$ws.Range("I136").Value2 = 8774518
$ws.Range("J136").Value2 = 5563.25
$ws.Range("K136").Value2 = 26323554
$ws.Range("L136").Value2 = 16689.75
$ws.Range("M136").Value2 = -26321004
$ws.Range("N136").Value2 = -21789.75
